# Daily attendance processing - reorder "Recorded By" author lists in column G
# so that the "System" placeholder entry is listed after the real user entry
# instead of before it (only for the specific recorder combinations that were
# corrected in this pass).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G
    $current = $cell.Value2

    if ($current -eq "System, dnasr281@gmail.com") {
        $cell.Value2 = "dnasr281@gmail.com, System"
    }
    elseif ($current -eq "System, system, backup@backdoor.com") {
        $cell.Value2 = "System, backup@backdoor.com, system"
    }
}
